$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.238.44"
$ws.Range("E2").Value = "  -6.11%  "
$ws.Range("D3").Value = "2.478.99"
$ws.Range("E3").Value = "  -7.79%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'541.67"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").Value = "'148.17"
$ws.Range("E6").Value = "  -6.92%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.572"
$ws.Range("E8").Value = "  -3.51%  "
$ws.Range("D9").Value = "2.476.53"
$ws.Range("E9").Value = "  -8.02%  "
$ws.Range("D10").Value = "'0.0995"
$ws.Range("E10").Value = "  -6.39%  "
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("D12").Value = "'5.34"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "2.903.38"
$ws.Range("E14").Value = "  -8.19%  "
$ws.Range("D15").Value = "'24.18"
$ws.Range("E15").Value = "  -9.05%  "
$ws.Range("D16").Value = "59.093.47"
$ws.Range("E16").Value = "  -6.18%  "
$ws.Range("E17").Value = "  -6.47%  "
$ws.Range("D18").Value = "2.523.58"
$ws.Range("E18").Value = "  -6.13%  "
$ws.Range("D19").Value = "'11.17"
$ws.Range("E19").Value = "  -7.16%  "
$ws.Range("E20").Value = "  -5.82%  "
$ws.Range("D21").Value = "'325.37"
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("D22").Value = "'0.968"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("E23").Value = "  -8.74%  "
$ws.Range("D24").Value = "'0.462"
$ws.Range("E24").Value = "  -9.94%  "
$ws.Range("D25").Value = "'60.79"
$ws.Range("E25").Value = "  -4.16%  "
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("D27").Value = "'0.980"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  -6.46%  "
$ws.Range("E29").Value = "  -9.77%  "
$ws.Range("E30").Value = "  -6.52%  "
$ws.Range("D31").Value = "'0.0₃0776"
$ws.Range("E31").Value = "  -10.55%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'158.01"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -8.37%  "
$ws.Range("E36").Value = "  -7.86%  "
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("D38").Value = "'1.74"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'5.92"
$ws.Range("E39").Value = "  -8.21%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'319.94"
$ws.Range("E40").Value = "  -10.94%  "
$ws.Range("D41").Value = "'36.61"
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("D42").Value = "'0.837"
$ws.Range("E42").Value = "  -13.09%  "
$ws.Range("D43").Value = "'3.71"
$ws.Range("E43").Value = "  -7.88%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'10.73"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("D47").Value = "'0.0945"
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("D48").Value = "'0.0527"
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("D49").Value = "'19.12"
$ws.Range("E49").Value = "  -9.67%  "
$ws.Range("E50").Value = "  -5.80%  "
$ws.Range("D51").Value = "'18.51"
$ws.Range("E51").Value = "  -9.63%  "
